# Update cryptocurrency price/volume data (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column price cells are plain text (e.g. "37.865.83", "0.380") that must
# stay text with exact digits -- temporarily force Text format so Excel
# does not reinterpret them as numbers (and drop trailing zeros), then
# restore the default "Normal" style so formatting is unchanged overall.

$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = "37.865.83"
$dCell.Style = "Normal"
$ws.Range("E2").Value = "  +6.30%  "

$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = "2.057.10"
$dCell.Style = "Normal"
$ws.Range("E3").Value = "  +3.78%  "

$ws.Range("E4").Value = "  -0.04%  "

$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = "253.21"
$dCell.Style = "Normal"
$ws.Range("E5").Value = "  +4.44%  "

$dCell = $ws.Range("D6")
$dCell.NumberFormat = "@"
$dCell.Value = "0.654"
$dCell.Style = "Normal"
$ws.Range("E6").Value = "  +2.52%  "

$dCell = $ws.Range("D7")
$dCell.NumberFormat = "@"
$dCell.Value = "65.46"
$dCell.Style = "Normal"
$ws.Range("E7").Value = "  +14.59%  "

$ws.Range("E8").Value = "  -0.01%  "

$dCell = $ws.Range("D9")
$dCell.NumberFormat = "@"
$dCell.Value = "60.88"
$dCell.Style = "Normal"
$ws.Range("E9").Value = "  +1.99%  "

$dCell = $ws.Range("D10")
$dCell.NumberFormat = "@"
$dCell.Value = "0.380"
$dCell.Style = "Normal"
$ws.Range("E10").Value = "  +5.63%  "

$dCell = $ws.Range("D11")
$dCell.NumberFormat = "@"
$dCell.Value = "0.0766"
$dCell.Style = "Normal"
$ws.Range("E11").Value = "  +4.95%  "

$ws.Range("E12").Value = "  +1.80%  "

$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value = "0.923"
$dCell.Style = "Normal"
$ws.Range("E13").Value = "  +0.06%  "

$dCell = $ws.Range("D14")
$dCell.NumberFormat = "@"
$dCell.Value = "15.07"
$dCell.Style = "Normal"
$ws.Range("E14").Value = "  +7.43%  "

$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = "2.358.57"
$dCell.Style = "Normal"
$ws.Range("E15").Value = "  +3.76%  "

$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value = "20.73"
$dCell.Style = "Normal"
$ws.Range("E16").Value = "  +20.56%  "

$dCell = $ws.Range("D17")
$dCell.NumberFormat = "@"
$dCell.Value = "5.57"
$dCell.Style = "Normal"
$ws.Range("E17").Value = "  +6.46%  "

$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"
$dCell.Value = "2.063.84"
$dCell.Style = "Normal"
$ws.Range("E18").Value = "  +4.15%  "

$dCell = $ws.Range("D19")
$dCell.NumberFormat = "@"
$dCell.Value = "37.760.08"
$dCell.Style = "Normal"
$ws.Range("E19").Value = "  +6.43%  "

$dCell = $ws.Range("D20")
$dCell.NumberFormat = "@"
$dCell.Value = "74.08"
$dCell.Style = "Normal"
$ws.Range("E20").Value = "  +4.82%  "

$ws.Range("E21").Value = "  +5.28%  "

$dCell = $ws.Range("D22")
$dCell.NumberFormat = "@"
$dCell.Value = "5.37"
$dCell.Style = "Normal"
$ws.Range("E22").Value = "  +6.15%  "

$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = "239.49"
$dCell.Style = "Normal"
$ws.Range("E23").Value = "  +2.60%  "

$ws.Range("E24").Value = "  +15.61%  "

$ws.Range("E25").Value = "  +0.02%  "

$ws.Range("E26").Value = "  +5.38%  "

$dCell = $ws.Range("D27")
$dCell.NumberFormat = "@"
$dCell.Value = "9.66"
$dCell.Style = "Normal"
$ws.Range("E27").Value = "  +5.88%  "

$dCell = $ws.Range("D28")
$dCell.NumberFormat = "@"
$dCell.Value = "160.17"
$dCell.Style = "Normal"
$ws.Range("E28").Value = "  -2.05%  "

$dCell = $ws.Range("D29")
$dCell.NumberFormat = "@"
$dCell.Value = "20.07"
$dCell.Style = "Normal"
$ws.Range("E29").Value = "  +3.35%  "

$dCell = $ws.Range("D30")
$dCell.NumberFormat = "@"
$dCell.Value = "0.116"
$dCell.Style = "Normal"
$ws.Range("E30").Value = "  +29.00%  "

$ws.Range("E31").Value = "  +2.56%  "

$dCell = $ws.Range("D32")
$dCell.NumberFormat = "@"
$dCell.Value = "5.23"
$dCell.Style = "Normal"
$ws.Range("E32").Value = "  +9.01%  "

$ws.Range("E33").Value = "  +6.82%  "

$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"
$dCell.Value = "4.74"
$dCell.Style = "Normal"
$ws.Range("E34").Value = "  +11.06%  "

$dCell = $ws.Range("D35")
$dCell.NumberFormat = "@"
$dCell.Value = "0.0620"
$dCell.Style = "Normal"
$ws.Range("E35").Value = "  +5.12%  "

$dCell = $ws.Range("D36")
$dCell.NumberFormat = "@"
$dCell.Value = "2.44"
$dCell.Style = "Normal"
$ws.Range("E36").Value = "  +3.16%  "

$ws.Range("E37").Value = "  +3.94%  "

$ws.Range("E38").Value = "  -0.11%  "

$dCell = $ws.Range("D39")
$dCell.NumberFormat = "@"
$dCell.Value = "6.12"
$dCell.Style = "Normal"
$ws.Range("E39").Value = "  +24.40%  "

$ws.Range("E40").Value = "  +17.08%  "

$dCell = $ws.Range("D41")
$dCell.NumberFormat = "@"
$dCell.Value = "2.83"
$dCell.Style = "Normal"
$ws.Range("E41").Value = "  +25.67%  "

$dCell = $ws.Range("D42")
$dCell.NumberFormat = "@"
$dCell.Value = "1.24"
$dCell.Style = "Normal"
$ws.Range("E42").Value = "  +4.47%  "

$ws.Range("E43").Value = "  +4.91%  "

$ws.Range("E44").Value = "  +3.63%  "

$ws.Range("E45").Value = "  +5.90%  "

$dCell = $ws.Range("D46")
$dCell.NumberFormat = "@"
$dCell.Value = "17.03"
$dCell.Style = "Normal"
$ws.Range("E46").Value = "  +10.15%  "

$dCell = $ws.Range("D47")
$dCell.NumberFormat = "@"
$dCell.Value = "8.00"
$dCell.Style = "Normal"
$ws.Range("E47").Value = "  +7.56%  "

$dCell = $ws.Range("D48")
$dCell.NumberFormat = "@"
$dCell.Value = "95.38"
$dCell.Style = "Normal"
$ws.Range("E48").Value = "  +4.75%  "

$dCell = $ws.Range("D49")
$dCell.NumberFormat = "@"
$dCell.Value = "1.415.53"
$dCell.Style = "Normal"
$ws.Range("E49").Value = "  +3.00%  "

$ws.Range("E50").Value = "  +2.47%  "

$dCell = $ws.Range("D51")
$dCell.NumberFormat = "@"
$dCell.Value = "47.38"
$dCell.Style = "Normal"
$ws.Range("E51").Value = "  +3.57%  "
